$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.310.67'
$ws.Range("E2").Value = '  -2.00%  '

$ws.Range("D3").Value = '3.689.11'
$ws.Range("E3").Value = '  -3.00%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''681.65'
$ws.Range("E5").Value = '  -3.77%  '

$ws.Range("D6").Value = '''162.84'
$ws.Range("E6").Value = '  -4.32%  '

$ws.Range("D7").Value = '3.685.41'
$ws.Range("E7").Value = '  -3.09%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.500'
$ws.Range("E9").Value = '  -4.03%  '

$ws.Range("E10").Value = '  -6.97%  '

$ws.Range("D11").Value = '''7.29'
$ws.Range("E11").Value = '  -1.53%  '

$ws.Range("D12").Value = '''0.451'
$ws.Range("E12").Value = '  -1.05%  '

$ws.Range("D13").Value = '''0.0000238'
$ws.Range("E13").Value = '  -5.80%  '

$ws.Range("D14").Value = '''33.66'
$ws.Range("E14").Value = '  -6.62%  '

$ws.Range("D15").Value = '4.309.81'
$ws.Range("E15").Value = '  -3.08%  '

$ws.Range("D16").Value = '3.687.04'
$ws.Range("E16").Value = '  -3.46%  '

$ws.Range("D17").Value = '69.346.07'
$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("E18").Value = '  -1.72%  '

$ws.Range("E19").Value = '  -5.91%  '

$ws.Range("D20").Value = '''6.66'
$ws.Range("E20").Value = '  -6.40%  '

$ws.Range("D21").Value = '''482.75'
$ws.Range("E21").Value = '  -2.55%  '

$ws.Range("D22").Value = '''9.79'
$ws.Range("E22").Value = '  -7.88%  '

$ws.Range("D23").Value = '''0.668'
$ws.Range("E23").Value = '  -8.50%  '

$ws.Range("D24").Value = '''79.99'
$ws.Range("E24").Value = '  -5.34%  '

$ws.Range("D25").Value = '3.834.33'
$ws.Range("E25").Value = '  -3.07%  '

$ws.Range("E26").Value = '  -10.31%  '

$ws.Range("D27").Value = '''11.55'
$ws.Range("E27").Value = '  -4.41%  '

$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").Value = '''9.63'
$ws.Range("E29").Value = '  -7.74%  '

$ws.Range("D30").Value = '''1.83'
$ws.Range("E30").Value = '  -10.32%  '

$ws.Range("E31").Value = '  -10.86%  '

$ws.Range("D32").Value = '''2.13'
$ws.Range("E32").Value = '  -4.86%  '

$ws.Range("D33").Value = '''6.84'
$ws.Range("E33").Value = '  -6.38%  '

$ws.Range("D34").Value = '''27.09'
$ws.Range("E34").Value = '  -6.57%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("E36").Value = '  -4.46%  '

$ws.Range("D37").Value = '3.651.30'
$ws.Range("E37").Value = '  -3.28%  '

$ws.Range("E38").Value = '  -5.92%  '

$ws.Range("D39").Value = '''6.10'
$ws.Range("E39").Value = '  +2.93%  '

$ws.Range("D40").Value = '''0.0947'
$ws.Range("E40").Value = '  -6.66%  '

$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("D42").Value = '''2.18'
$ws.Range("E42").Value = '  -5.47%  '

$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("D44").Value = '''0.959'
$ws.Range("E44").Value = '  -7.68%  '

$ws.Range("D45").Value = '''158.01'
$ws.Range("E45").Value = '  -4.23%  '

$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("D47").Value = '''2.86'
$ws.Range("E47").Value = '  -12.05%  '

$ws.Range("E48").Value = '  -12.15%  '

$ws.Range("D49").Value = '''1.31'
$ws.Range("E49").Value = '  -3.52%  '

$ws.Range("D50").Value = '''389.44'
$ws.Range("E50").Value = '  -8.46%  '

$ws.Range("D51").Value = '''8.11'
$ws.Range("E51").Value = '  -5.72%  '

